$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the answer text in B15 (Serdar Altay research-assistant question)
# so it now repeats the project name ("Serdar Altay'ın Trade Policy Review of
# Turkey isimli projesi için ...") instead of "Prof. Serdar Altay için ...".
$ws.Range("B15").Value = "Serdar Altay'ın Trade Policy Review of Turkey isimli projesi için yaptığımız çalışma veri analizi, toplaması ve editöryel destek sayesinde hocamız çalışmasını tamamlamış ve bu çalışma The World Economy isimli dergide yayınlanmıştır."

# The longer text now needs a taller row to display it (wrap text is already
# enabled via the existing cell style).
$ws.Rows.Item(15).RowHeight = 85

# Scroll the view so row 9 is at the top and select B15, matching where the
# user was working when they made this edit.
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("B15").Select()
